$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.923.56'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '1.815.22'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.80'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4687'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3690'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07366'
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8727'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.38'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '1.840.61'
$ws.Range("E12").Value = '  +4.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.375'
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07085'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.517'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.89'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008709'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.71'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '26.953.97'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.325'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.62'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.25%  '
$ws.Range("D24").Value = '2.020.77'
$ws.Range("E24").Value = '  +2.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.889'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.64'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.179'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.36'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.320'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.06'
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08949'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7654'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.166'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.502'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.921'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("E37").Value = '  -2.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01960'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05288'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.80%  '
$ws.Range("E40").Value = '  +2.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.258'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5344'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.324'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.449'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4924'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.42'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.669'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.06'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("E51").Value = '  -0.20%  '
